# Auto update Excel log
# Appends 5 new PRESENCE_DETECTED sensor readings (rows 36-40) to the
# "mmWave" sheet, matching the source log export. The Date/Timestamp/Hour
# columns look like dates/times to Excel's auto-detection, so each of
# those cells is temporarily switched to Text format before the value is
# written (and the style is reset back to "Normal" afterward) to make
# sure they land as literal strings instead of being coerced into date
# serial numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mmWave")

$newRows = @(
    @("2026-01-31", "21:45:13", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:45:13", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:45:14", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:45:15", "21:00", "Living Room", "PRESENCE_DETECTED", "Active"),
    @("2026-01-31", "21:45:26", "21:00", "Living Room", "PRESENCE_DETECTED", "Active")
)

$startRow = 36
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]

    for ($col = 1; $col -le 3; $col++) {
        $cell = $ws.Cells.Item($r, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $values[$col - 1]
        $cell.Style = "Normal"
    }

    for ($col = 4; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $values[$col - 1]
    }
}
